# "Daniels assesments merged in to my structure"
# Fill in Daniel's peer-assessment entries on the
# "Peer  and self assessment" sheet (rows 8 and 20),
# which were previously left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")

# Criterion 1 Online collaboration -> Daniel
$ws.Range("B8").Value = "Good"
$ws.Range("C8").Value = "Leading group meetings, fast response times and very active on discord."

# Criterion 1 International Collaboration -> Daniel
$ws.Range("B20").Value = "Excellent"
$ws.Range("C20").Value = "Active collaborator, motivated"
